# Generate Report for Handoff
#
# A new handoff was generated for the files that are not already fully
# "Handed back: in sync with en-US" and are not the ignored
# ".localization-config" entry. For those rows (7-16 on each language
# status sheet) the "Latest Handoff Datetime" (column D) is stamped with
# the new handoff timestamp for that language.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnHandoffTime = "2016-03-09 12:24:34"
for ($row = 7; $row -le 16; $row++) {
    $zhcn.Range("D$row").Value = $zhcnHandoffTime
}

$dede = $wb.Worksheets.Item("de-de")
$dedeHandoffTime = "2016-03-09 12:24:40"
for ($row = 7; $row -le 16; $row++) {
    $dede.Range("D$row").Value = $dedeHandoffTime
}
